$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.895.98'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '1.821.09'
$ws.Range('E4').Value = '  -0.86%  '
$ws.Range('D5').Value = "'309.76"
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').Value = "'0.4623"
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').Value = "'0.3685"
$ws.Range('E8').Value = '  +1.36%  '
$ws.Range('D9').Value = "'0.07308"
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('D10').Value = "'0.8734"
$ws.Range('E10').Value = '  +1.06%  '
$ws.Range('D11').Value = "'0.07843"
$ws.Range('E11').Value = '  +3.19%  '
$ws.Range('D12').Value = "'19.57"
$ws.Range('E12').Value = '  -1.03%  '
$ws.Range('D13').Value = '1.791.65'
$ws.Range('E13').Value = '  -4.00%  '
$ws.Range('D14').Value = "'5.320"
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('D15').Value = "'6.529"
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').Value = "'91.04"
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('D18').Value = "'0.000008831"
$ws.Range('E18').Value = '  +2.52%  '
$ws.Range('D19').Value = "'1.000"
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('D20').Value = "'14.70"
$ws.Range('E20').Value = '  +1.71%  '
$ws.Range('D21').Value = '26.919.03'
$ws.Range('E21').Value = '  -1.62%  '
$ws.Range('D22').Value = "'5.090"
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('D24').Value = '2.011.63'
$ws.Range('E24').Value = '  -4.58%  '
$ws.Range('D25').Value = "'151.79"
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = "'1.851"
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').Value = "'18.28"
$ws.Range('E27').Value = '  +0.49%  '
$ws.Range('D28').Value = "'2.028"
$ws.Range('E28').Value = '  -2.71%  '
$ws.Range('D29').Value = "'5.084"
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').Value = "'114.93"
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('D31').Value = "'0.08834"
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').Value = "'2.954"
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').Value = "'0.7296"
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('D34').Value = "'4.420"
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = "'1.128"
$ws.Range('E35').Value = '  -0.96%  '
$ws.Range('D36').Value = "'2.458"
$ws.Range('E36').Value = '  -1.33%  '
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('D38').Value = "'0.01933"
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').Value = "'0.05210"
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('D40').Value = "'2.951"
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').Value = "'7.049"
$ws.Range('E41').Value = '  -1.34%  '
$ws.Range('D42').Value = "'0.5112"
$ws.Range('E42').Value = '  -1.51%  '
$ws.Range('D43').Value = "'0.1618"
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').Value = "'0.4819"
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = "'8.120"
$ws.Range('E45').Value = '  -1.57%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = "'1.002"
$ws.Range('E46').Value = '  -0.76%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'10.14"
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('D48').Value = "'101.49"
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('D49').Value = "'1.615"
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('D50').Value = "'0.06190"
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('D51').Value = "'64.27"
$ws.Range('E51').Value = '  -0.51%  '
